$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-28: increment date by 1 day (45458 -> 45459)
for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 3).Value = 45459
}
